# "new - 更新 3.12 中签率" — fill in the winning-rate (中签率) and its
# computed total-chance figures for the five IPOs that were still blank
# on the 2015-03 sheet (rows 8-12), mirroring the pattern already used
# for rows 2 and 7: column I holds the disclosed winning rate, column J
# computes I*G (shares allotted * rate).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2015-03")

# Carry over the exact number format / font used by the existing
# "中签率" entries (I7 / J7) onto the newly-filled rows so the new
# cells render identically (red percentage text, boxed border).
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I8:I12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("J7").Copy() | Out-Null
$ws.Range("J8:J12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 暴风科技 SZ-300431
$ws.Range("I8").Value = 3.4E-3
# 龙津药业 SZ-002750
$ws.Range("I9").Value = 6.4E-3
# 埃斯顿 SZ-002747
$ws.Range("I10").Value = 3.9E-3
# 国光股份 SZ-002749
$ws.Range("I11").Value = 6.8E-3
# 强力新材 SZ-300429
$ws.Range("I12").Value = 7.4E-3

# J = 中签率(I) * 配号数(G)
$ws.Range("J8:J12").FormulaR1C1 = "=RC[-1]*RC[-3]"

# The workbook was last saved with L13 selected.
$ws.Range("L13").Select() | Out-Null
